$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1266
$ws.Range("I100").Value = 1209.0769
$ws.Range("J100").Value = 2006
$ws.Range("K100").Value = 1209.0769
$ws.Range("L100").Value = 2006
$ws.Range("M100").Value = -668.0769
$ws.Range("N100").Value = -3088
$ws.Range("H107").Value = 1147.3
$ws.Range("I107").Value = 1147.3
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1147.3
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 772.7
$ws.Range("N107").ClearContents()
$ws.Range("H113").Value = 166690500
$ws.Range("I113").Value = 200003580
$ws.Range("J113").Value = 125000
$ws.Range("K113").Value = 200003580
$ws.Range("L113").Value = 125000
$ws.Range("M113").Value = -200000326
$ws.Range("N113").Value = -131508
$ws.Range("H131").Value = 3472.0588
$ws.Range("I131").Value = 2439.0625
$ws.Range("J131").Value = 20000
$ws.Range("K131").Value = 7317.1875
$ws.Range("L131").Value = 60000
$ws.Range("M131").Value = -2277.1875
$ws.Range("N131").Value = -70080
$ws.Range("H132").Value = 9633.879999999999
$ws.Range("I132").Value = 1783.238
$ws.Range("J132").Value = 50849.75
$ws.Range("K132").Value = 5349.714
$ws.Range("L132").Value = 152549.25
$ws.Range("M132").Value = -2819.714
$ws.Range("N132").Value = -157609.25
$ws.Range("H137").Value = 4195.409
$ws.Range("I137").Value = 3266.3333
$ws.Range("J137").Value = 4342.1055
$ws.Range("K137").Value = 9798.999899999999
$ws.Range("L137").Value = 13026.3165
$ws.Range("M137").Value = -7248.999899999999
$ws.Range("N137").Value = -18126.3165
$ws.Range("H138").Value = 3476.0852
$ws.Range("I138").Value = 1028.7407
$ws.Range("J138").Value = 6780
$ws.Range("K138").Value = 3086.2221
$ws.Range("L138").Value = 20340
$ws.Range("M138").Value = 2053.7779
$ws.Range("N138").Value = -30620
$ws.Range("H141").Value = 3190.258
$ws.Range("I141").Value = 3172.6072
$ws.Range("J141").Value = 3355
$ws.Range("K141").Value = 9517.821599999999
$ws.Range("L141").Value = 10065
$ws.Range("M141").Value = -4337.821599999999
$ws.Range("N141").Value = -20425

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3816.8823
$ws.Range("I2").Value = 2760.6924
$ws.Range("K2").Value = 2760.6924
$ws.Range("M2").Value = -2647.6924
$ws.Range("H32").Value = 6569.5884
$ws.Range("I32").Value = 6026.231
$ws.Range("K32").Value = 6026.231
$ws.Range("M32").Value = -5739.231
$ws.Range("H61").Value = 4091.0876
$ws.Range("I61").Value = 4277.34
$ws.Range("J61").Value = 1623.25
$ws.Range("K61").Value = 4277.34
$ws.Range("L61").Value = 1623.25
$ws.Range("M61").Value = -4065.34
$ws.Range("N61").Value = -2047.25
$ws.Range("H74").Value = 3403.5957
$ws.Range("I74").Value = 2346.4324
$ws.Range("J74").Value = 7315.1
$ws.Range("K74").Value = 2346.4324
$ws.Range("L74").Value = 7315.1
$ws.Range("M74").Value = -1472.4324
$ws.Range("N74").Value = -9063.1
$ws.Range("H77").Value = 3403.5957
$ws.Range("I77").Value = 2346.4324
$ws.Range("J77").Value = 7315.1
$ws.Range("K77").Value = 11732.162
$ws.Range("L77").Value = 36575.5
$ws.Range("M77").Value = -7364.162
$ws.Range("N77").Value = -45311.5
$ws.Range("H97").Value = 1556.7858
$ws.Range("I97").Value = 1285.762
$ws.Range("J97").Value = 2369.8572
$ws.Range("K97").Value = 1285.762
$ws.Range("L97").Value = 2369.8572
$ws.Range("M97").Value = -789.7619999999999
$ws.Range("N97").Value = -3361.8572
$ws.Range("H110").Value = 4549.5
$ws.Range("I110").Value = 3974.9583
$ws.Range("J110").Value = 7996.75
$ws.Range("K110").Value = 3974.9583
$ws.Range("L110").Value = 7996.75
$ws.Range("M110").Value = -1929.9583
$ws.Range("N110").Value = -12086.75
$ws.Range("H116").Value = 3816.8823
$ws.Range("I116").Value = 2760.6924
$ws.Range("K116").Value = 2760.6924
$ws.Range("M116").Value = -466.6923999999999
$ws.Range("H136").Value = 4091.0876
$ws.Range("I136").Value = 4277.34
$ws.Range("J136").Value = 1623.25
$ws.Range("K136").Value = 12832.02
$ws.Range("L136").Value = 4869.75
$ws.Range("M136").Value = -10282.02
$ws.Range("N136").Value = -9969.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3816.8823
$ws.Range("I3").Value = 2760.6924
$ws.Range("K3").Value = 2760.6924
$ws.Range("M3").Value = -2646.6924
$ws.Range("H86").Value = 1296.3438
$ws.Range("I86").Value = 1206.96
$ws.Range("J86").Value = 1615.5714
$ws.Range("K86").Value = 1206.96
$ws.Range("L86").Value = 1615.5714
$ws.Range("M86").Value = -83.96000000000004
$ws.Range("N86").Value = -3861.5714
$ws.Range("H89").Value = 1296.3438
$ws.Range("I89").Value = 1206.96
$ws.Range("J89").Value = 1615.5714
$ws.Range("K89").Value = 6034.8
$ws.Range("L89").Value = 8077.857
$ws.Range("M89").Value = -418.8000000000002
$ws.Range("N89").Value = -19309.857
$ws.Range("H134").Value = 3518.0676
$ws.Range("I134").Value = 3415.9033
$ws.Range("J134").Value = 4045.9167
$ws.Range("K134").Value = 10247.7099
$ws.Range("L134").Value = 12137.7501
$ws.Range("M134").Value = -7712.7099
$ws.Range("N134").Value = -17207.7501

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3399.5
$ws.Range("I16").Value = 2666.111
$ws.Range("J16").Value = 10000
$ws.Range("K16").Value = 2666.111
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = -2379.111
$ws.Range("N16").Value = -10574
$ws.Range("H31").Value = 26320364
$ws.Range("I31").Value = 50003068
$ws.Range("J31").Value = 6246.722
$ws.Range("K31").Value = 50003068
$ws.Range("L31").Value = 6246.722
$ws.Range("M31").Value = -50002773
$ws.Range("N31").Value = -6836.722
$ws.Range("H34").Value = 26320364
$ws.Range("I34").Value = 50003068
$ws.Range("J34").Value = 6246.722
$ws.Range("K34").Value = 50003068
$ws.Range("L34").Value = 6246.722
$ws.Range("M34").Value = -50002866
$ws.Range("N34").Value = -6650.722
$ws.Range("H58").Value = 13368.125
$ws.Range("I58").Value = 9195
$ws.Range("J58").Value = 14759.167
$ws.Range("K58").Value = 9195
$ws.Range("L58").Value = 14759.167
$ws.Range("M58").Value = -8992
$ws.Range("N58").Value = -15165.167
$ws.Range("H113").Value = 3399.5
$ws.Range("I113").Value = 2666.111
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 2666.111
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = -496.1109999999999
$ws.Range("N113").Value = -14340
$ws.Range("H115").Value = 40290
$ws.Range("J115").Value = 40290
$ws.Range("L115").Value = 40290
$ws.Range("N115").Value = -42640
$ws.Range("H132").Value = 2873.4897
$ws.Range("I132").Value = 1718.85
$ws.Range("J132").Value = 8005.222
$ws.Range("K132").Value = 5156.549999999999
$ws.Range("L132").Value = 24015.666
$ws.Range("M132").Value = -2626.549999999999
$ws.Range("N132").Value = -29075.666
$ws.Range("H134").Value = 7108.9116
$ws.Range("I134").Value = 5471.5356
$ws.Range("J134").Value = 14750
$ws.Range("K134").Value = 16414.6068
$ws.Range("L134").Value = 44250
$ws.Range("M134").Value = -13879.6068
$ws.Range("N134").Value = -49320
$ws.Range("H136").Value = 13368.125
$ws.Range("I136").Value = 9195
$ws.Range("J136").Value = 14759.167
$ws.Range("K136").Value = 27585
$ws.Range("L136").Value = 44277.501
$ws.Range("M136").Value = -25035
$ws.Range("N136").Value = -49377.501

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3053.8333
$ws.Range("I113").Value = 3053.8333
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3053.8333
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -883.8332999999998
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 7249.75
$ws.Range("I132").Value = 7333
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 21999
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -19469
$ws.Range("N132").Value = -26060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3285.6592
$ws.Range("I136").Value = 3491.725
$ws.Range("J136").Value = 1225
$ws.Range("K136").Value = 10475.175
$ws.Range("L136").Value = 3675
$ws.Range("M136").Value = -7925.174999999999
$ws.Range("N136").Value = -8775

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2529.4814
$ws.Range("I136").Value = 2014.5952
$ws.Range("J136").Value = 4331.5835
$ws.Range("K136").Value = 6043.7856
$ws.Range("L136").Value = 12994.7505
$ws.Range("M136").Value = -3493.7856
$ws.Range("N136").Value = -18094.7505

